$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("foot_stats")

# Columns J (10) and K (11) were re-purposed: J used to hold the scorer's
# name (text) with K holding the goal count (number); now K holds the name
# and J holds the number. Swap every data row's J/K content accordingly.
for ($r = 2; $r -le 99; $r++) {
    $jCell = $ws.Cells.Item($r, 10)
    $kCell = $ws.Cells.Item($r, 11)
    $jVal = $jCell.Value()
    $kVal = $kCell.Value()
    $jCell.Value = $kVal
    $kCell.Value = $jVal
}

# Header row: J1 becomes the new "QG" label, K1 takes over the old
# "Top Team Scorer" label (previously in J1).
$ws.Cells.Item(1, 10).Value = "QG"
$ws.Cells.Item(1, 11).Value = "Top Team Scorer"

# Cosmetic: selection moves to N6, and J/K get explicit column widths.
[void]$ws.Range("N6").Select()
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 16.5
